$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header summary cells: mora total, worker count, period count
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 7028138
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 80

# ---------------------------------------------------------------------------
# 2. Capture the existing "Periodo Mora" column (rows 16-94, CESAR ENRIQUE's
#    79 periods, ascending 1812..2506) before we touch the layout.
# ---------------------------------------------------------------------------
$oldPeriods = @()
for ($r = 16; $r -le 94; $r++) {
    $oldPeriods += $ws.Cells.Item($r, 5).Value()
}

# ---------------------------------------------------------------------------
# 3. Insert three new rows right after the existing data block (was rows
#    16-94, ending at the "last row" thick-border style). This pushes the
#    closing signature block (old rows 99-100) down to 102-103.
# ---------------------------------------------------------------------------
$ws.Range("95:97").Insert()

# Propagate formatting: rows 95 & 96 get the regular interior-row style
# (copied from row 93); row 97 becomes the new thick-bottom-border closing
# row (copied from the old row 94); row 94 itself is downgraded from the
# closing style back to the regular interior style.
$ws.Range("B93:J93").Copy($ws.Range("B95"))
$ws.Range("B93:J93").Copy($ws.Range("B96"))
$ws.Range("B94:J94").Copy($ws.Range("B97"))
$ws.Range("B93:J93").Copy($ws.Range("B94"))

# ---------------------------------------------------------------------------
# 4. Rewrite the "Periodo Mora" column for CESAR ENRIQUE GALVAN GUERRA:
#    newest period "2507" first, then the previous periods in reverse
#    (descending) order, i.e. 2506, 2505, ... down to 1812 -- 80 rows total
#    (16 through 95).
# ---------------------------------------------------------------------------
$ws.Cells.Item(16, 5).Value = "2507"
for ($i = 0; $i -lt $oldPeriods.Count; $i++) {
    $r = 17 + $i
    $period = $oldPeriods[$oldPeriods.Count - 1 - $i]
    $ws.Cells.Item($r, 5).Value = $period
}

# ---------------------------------------------------------------------------
# 5. New worker block: LEONARDO ENRIQUE CARCAMO VEGA, 2 overdue periods.
# ---------------------------------------------------------------------------
$ws.Cells.Item(96, 2).Value = "CC"
$ws.Cells.Item(96, 3).Value = "73236799"
$ws.Cells.Item(96, 4).Value = "LEONARDO ENRIQUE CARCAMO VEGA"
$ws.Cells.Item(96, 5).Value = "2507"
$ws.Cells.Item(96, 6).Value = 6150
$ws.Cells.Item(96, 7).Value = 4612606

$ws.Cells.Item(97, 2).Value = "CC"
$ws.Cells.Item(97, 3).Value = "73236799"
$ws.Cells.Item(97, 4).Value = "LEONARDO ENRIQUE CARCAMO VEGA"
$ws.Cells.Item(97, 5).Value = "2505"
$ws.Cells.Item(97, 6).Value = 5748
$ws.Cells.Item(97, 7).Value = 4612606
